$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("D3").Value = "Task"
$ws.Range("E3").Value = "Hours Spent"

# Task rows
$ws.Range("D4").Value = "Pileup Parser"
$ws.Range("E4").Value = 6

$ws.Range("D5").Value = "SVG Creator"
$ws.Range("E5").Value = 3

$ws.Range("D6").Value = "Flow Diagram"
$ws.Range("E6").Value = 2

# Update selection to match the new used range
$ws.Range("D3:E6").Select()
